# excel & widget refacto nature0 groupe0 & limite0
#
# 1) "lines" sheet: fix the P0-E1 row (dist/durite), center B:C columns,
#    move the selection.
# 2) add a new "slot" sheet (capteur / group / nature / limit).

$wb = $excel.ActiveWorkbook
$linesWs = $wb.Worksheets.Item("lines")

# --- correct the P0-E1 row values (row 12) ---
$linesWs.Range("C12").Value = 0.1
$linesWs.Range("D12").Value = 4

# --- center-align columns B and C (data rows only, header row stays as-is) ---
$linesWs.Range("B2:C19").HorizontalAlignment = -4108

# --- give columns B:C the same default width/style as column D ---
$linesWs.Range("B:C").EntireColumn.ColumnWidth = 9.140625

# --- move the active selection on the "lines" sheet ---
$linesWs.Activate() | Out-Null
$linesWs.Range("G8").Select() | Out-Null

# --- add the new "slot" sheet after "lines" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$slotWs = $wb.Worksheets.Add($null, $lastSheet)
$slotWs.Name = "slot"

# header row (note: group is written before capteur so the shared-string
# table is rebuilt in the same order as the source workbook)
$slotWs.Range("C1").Value = "group"
$slotWs.Range("B1").Value = "capteur"
$slotWs.Range("D1").Value = "nature"
$slotWs.Range("E1").Value = "limit"

# data rows
$slotWs.Range("A2").Value = 0
$slotWs.Range("B2").Value = "C0"
$slotWs.Range("C2").Value = 1
$slotWs.Range("D2").Value = "F"
$slotWs.Range("E2").Value = 1.7

$slotWs.Range("A3").Value = 1
$slotWs.Range("B3").Value = "C1"
$slotWs.Range("C3").Value = 1
$slotWs.Range("D3").Value = "R"
$slotWs.Range("E3").Value = 2

$slotWs.Range("A4").Value = 2
$slotWs.Range("B4").Value = "C2"
$slotWs.Range("C4").Value = 0
$slotWs.Range("D4").Value = "F"
$slotWs.Range("E4").Value = 2

$slotWs.Range("A5").Value = 3
$slotWs.Range("B5").Value = "C3"
$slotWs.Range("C5").Value = 0
$slotWs.Range("D5").Value = "F"
$slotWs.Range("E5").Value = 2

# formatting: A2:D5 + the row-1 headers (B1:D1 and E1) are centered; the
# E1 label is centered too, but the E2:E5 data cells stay General. A1 is
# never touched (there is no header above the "capteur" id column).
$slotWs.Range("A2:D5").HorizontalAlignment = -4108
$slotWs.Range("B1:E1").HorizontalAlignment = -4108
$slotWs.Range("A:D").EntireColumn.ColumnWidth = 9.140625

# the new sheet remembers F5 as its own selection, but stays in the
# background -- re-activate "lines" so it is still the visible tab
$slotWs.Range("F5").Select() | Out-Null
$linesWs.Activate() | Out-Null
